$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "alimentación" (data entry) for tracking row 2 (sheet row 7):
# dates for Publicación / Manuscrito / Formatos de recursos / Guión, and the
# "Entrega" observation note.
$ws.Range("B7").Value = "3/16/2015"
$ws.Range("C7").Value = "3/16/2015"
$ws.Range("D7").Value = "3/16/2015"
$ws.Range("E7").Value = "3/16/2015"
$ws.Range("F7").Value = "Completo después de corrección estilo."

# The row now holds wrapped text, so it grows to match the taller rows above.
$ws.Rows("7").RowHeight = 30.75

# Leave the selection on the "Grado" field, matching where the user was
# working.
$ws.Activate()
$ws.Range("C2:D2").Select()
